$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Formula = "=`"" + $text + "`""
    $rng.Copy($null)
    $rng.PasteSpecial(-4163, $null, $false, $false)
}

$ws.Range("D2").Value = "62.049.35"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.421.27"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "562.08"
$ws.Range("E5").Value = "  -0.04%  "
Set-TextValue "D6" "143.74"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "2.420.47"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  +0.27%  "
Set-TextValue "D12" "5.21"
$ws.Range("E12").Value = "  -3.20%  "
Set-TextValue "D13" "0.348"
$ws.Range("E13").Value = "  -1.38%  "
Set-TextValue "D14" "26.10"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "61.957.89"
$ws.Range("D18").Value = "2.412.04"
$ws.Range("E18").Value = "  -0.35%  "
Set-TextValue "D19" "11.23"
$ws.Range("E19").Value = "  -0.67%  "
Set-TextValue "D20" "323.05"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  -1.67%  "
Set-TextValue "D22" "6.82"
$ws.Range("E22").Value = "  +1.04%  "
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  -0.02%  "
Set-TextValue "D24" "67.43"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +1.28%  "
Set-TextValue "D26" "8.67"
$ws.Range("E26").Value = "  -3.02%  "
Set-TextValue "D27" "559.06"
$ws.Range("E27").Value = "  -4.45%  "
$ws.Range("D28").Value = "2.540.02"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  -1.34%  "
Set-TextValue "D31" "8.18"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  -4.86%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -4.48%  "
Set-TextValue "D40" "152.21"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -1.33%  "
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  +0.12%  "
Set-TextValue "D44" "2.26"
$ws.Range("E44").Value = "  -3.08%  "
Set-TextValue "D45" "147.08"
$ws.Range("E45").Value = "  -2.18%  "
Set-TextValue "D46" "3.64"
$ws.Range("E46").Value = "  -0.71%  "
Set-TextValue "D48" "19.94"
$ws.Range("E48").Value = "  -2.08%  "
Set-TextValue "D49" "0.594"
$ws.Range("E49").Value = "  -0.07%  "
Set-TextValue "D50" "0.0921"
$ws.Range("E50").Value = "  -0.35%  "
Set-TextValue "D51" "0.0228"
$ws.Range("E51").Value = "  -0.51%  "

# Rows 37/38: coin data swapped (NEARProtocol <-> PolygonEcosystemToken)
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D37" "4.74"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D38" "0.380"
$ws.Range("E38").Value = "  -1.04%  "
